# For each bus/device location on the "Bus" sheet, the phase-A row and the
# phase-C row were swapped (name in column A and angle in column E), while
# the other columns (Base Voltage, Initial Vmag, Unit, Type) stayed put
# since they were identical for the A/C pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bus")

# Row pairs (firstRow = was "_A"/0deg, secondRow = was "_C"/120deg)
$pairs = @(4,5), @(7,8), @(10,11), @(17,18), @(21,22), @(24,25), @(27,28), @(30,31), @(32,33), @(35,36), @(38,39), @(43,44), @(46,47), @(49,50), @(53,54)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $name1 = $ws.Cells.Item($r1, 1).Value2
    $angle1 = $ws.Cells.Item($r1, 5).Value2
    $name2 = $ws.Cells.Item($r2, 1).Value2
    $angle2 = $ws.Cells.Item($r2, 5).Value2

    $ws.Cells.Item($r1, 1).Value = $name2
    $ws.Cells.Item($r1, 5).Value = $angle2
    $ws.Cells.Item($r2, 1).Value = $name1
    $ws.Cells.Item($r2, 5).Value = $angle1
}
